$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (row 721), shifting all
# existing data rows down by two. This mirrors a new week of price
# observations being prepended to the dataset.
$ws.Rows("721:722").Insert()

# New row 721: Papa / Asterix (Provincia de Llanquihue)
$ws.Range("A721").Value = 4
$ws.Range("B721").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C721").Value = 'Los Lagos'
$ws.Range("D721").Value = 45223
$ws.Range("E721").Value = 10
$ws.Range("F721").Value = 100114001
$ws.Range("G721").Value = 'Papa'
$ws.Range("H721").Value = 'Asterix'
$ws.Range("I721").Value = '1a (guarda)'
$ws.Range("J721").Value = 300
$ws.Range("K721").Value = 29000
$ws.Range("L721").Value = 30000
$ws.Range("M721").Value = 29500
$ws.Range("N721").Value = '$/saco 25 kilos'
$ws.Range("O721").Value = 'Provincia de Llanquihue'
$ws.Range("P721").Value = 1180
$ws.Range("Q721").Value = 25
$ws.Range("R721").Value = 'Hortaliza'

# New row 722: Papa / Patagonia (Región de La Araucanía)
$ws.Range("A722").Value = 4
$ws.Range("B722").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C722").Value = 'Los Lagos'
$ws.Range("D722").Value = 45223
$ws.Range("E722").Value = 10
$ws.Range("F722").Value = 100114001
$ws.Range("G722").Value = 'Papa'
$ws.Range("H722").Value = 'Patagonia'
$ws.Range("I722").Value = '1a nueva(o)'
$ws.Range("J722").Value = 200
$ws.Range("K722").Value = 29000
$ws.Range("L722").Value = 30000
$ws.Range("M722").Value = 29500
$ws.Range("N722").Value = '$/saco 25 kilos'
$ws.Range("O722").Value = 'Región de La Araucanía'
$ws.Range("P722").Value = 1180
$ws.Range("Q722").Value = 25
$ws.Range("R722").Value = 'Hortaliza'

